# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The sheet currently holds player/team stats in columns A:AB (row 1 is the
# header row, rows 2:43 are the players). We append three new columns:
#   AC -> Wins
#   AD -> Losses
#   AE -> Ties
# and fill every player row with the team's 1992 season record (78-84-0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the formatting of an existing header cell (bold, centered, bordered)
# onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows ----------------------------------------------------------
# Every player on the roster shares the same team season record.
$wins = 78
$losses = 84
$ties = 0

$lastRow = 43
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins    # column AC
    $ws.Cells.Item($row, 30).Value = $losses  # column AD
    $ws.Cells.Item($row, 31).Value = $ties    # column AE
}
